# Rellenar la columna "uds. Objetivo semana pasada" (columna Q) para las
# filas correspondientes, segun los nuevos valores calculados.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_8")

$updates = @{
    2  = 1
    3  = 7
    4  = 1
    6  = 2
    8  = 3
    9  = 1
    11 = 3
    12 = 1
    16 = 1
    17 = 1
    18 = 1
    22 = 4
    23 = 2
    24 = 4
    26 = 7
    27 = 9
    28 = 4
    29 = 2
    32 = 1
    35 = 1
    36 = 2
    37 = 3
    39 = 5
    40 = 2
    44 = 2
    45 = 3
    47 = 1
    49 = 3
    50 = 3
    55 = 5
    57 = 8
    65 = 1
    67 = 2
    76 = 1
    77 = 34
    78 = 39
    79 = 3
    81 = 8
    82 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("Q$row").Value = $updates[$row]
}
